# Insert two new rows (265 and 266) above the current row 265, shifting the
# existing data (old rows 265-285) down to rows 267-287, then populate the
# two new rows with the latest weekly price entries.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 2 blank rows at 265:266 - this pushes old rows 265-285 to 267-287.
$ws.Rows("265:266").Insert()

# New row 265 - "Primera" quality entry for the latest week.
$ws.Range("A265").Value = 8
$ws.Range("B265").Value = "Terminal La Palmera de La Serena"
$ws.Range("C265").Value = "Coquimbo"
$ws.Range("D265").Value = 44826
$ws.Range("E265").Value = 4
$ws.Range("F265").Value = 100112021
$ws.Range("G265").Value = "Ají"
$ws.Range("H265").Value = "Inferno"
$ws.Range("I265").Value = "Primera"
$ws.Range("J265").Value = 600
$ws.Range("K265").Value = 21000
$ws.Range("L265").Value = 22000
$ws.Range("M265").Value = 21500
$ws.Range("N265").Value = '$/caja 10 kilos'
$ws.Range("O265").Value = "Región de Arica y Parinacota"
$ws.Range("P265").Value = 2150
$ws.Range("Q265").Value = 10
$ws.Range("R265").Value = "Hortaliza"

# New row 266 - "Segunda" quality entry for the latest week.
$ws.Range("A266").Value = 8
$ws.Range("B266").Value = "Terminal La Palmera de La Serena"
$ws.Range("C266").Value = "Coquimbo"
$ws.Range("D266").Value = 44826
$ws.Range("E266").Value = 4
$ws.Range("F266").Value = 100112021
$ws.Range("G266").Value = "Ají"
$ws.Range("H266").Value = "Inferno"
$ws.Range("I266").Value = "Segunda"
$ws.Range("J266").Value = 400
$ws.Range("K266").Value = 13000
$ws.Range("L266").Value = 14000
$ws.Range("M266").Value = 13500
$ws.Range("N266").Value = '$/caja 10 kilos'
$ws.Range("O266").Value = "Región de Arica y Parinacota"
$ws.Range("P266").Value = 1350
$ws.Range("Q266").Value = 10
$ws.Range("R266").Value = "Hortaliza"
